$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

# --- Add two more days (Z = 23-Jan-2024, AA = 24-Jan-2024) ---
# Copy the header cell (Y1) format into the new header cells, then set the date values
$ws.Range("Y1").Copy($ws.Range("Z1:AA1"))
$ws.Range("Z1").Value = 45314
$ws.Range("AA1").Value = 45315

# Copy the attendance-cell format (from column Y) into the new Z/AA columns for
# each student row, then fill in the attendance values for 23-Jan and 24-Jan.
$ws.Range("Y2").Copy($ws.Range("Z2"))
$ws.Range("Z2").Value = "Absent"
$ws.Range("Y2").Copy($ws.Range("AA2"))
$ws.Range("AA2").Value = "Present"

$ws.Range("Y3").Copy($ws.Range("Z3"))
$ws.Range("Z3").Value = "Absent"
$ws.Range("Y3").Copy($ws.Range("AA3"))
$ws.Range("AA3").Value = "Present"

$ws.Range("Y4").Copy($ws.Range("Z4"))
$ws.Range("Z4").Value = "Absent"
$ws.Range("Y4").Copy($ws.Range("AA4"))
$ws.Range("AA4").Value = "Present"

# --- Extend the "Present, Absent, Reason" dropdown validation to cover the new columns ---
$ws.Range("C2:Y4").Validation.Delete()
$ws.Range("C2:AA4").Validation.Add(3, 1, 1, """Present, Absent,Reason""")

# --- Update the view: scroll right so column Q is at the left edge, and move the
#     active selection down to AB10 (just past the new data) ---
$excel.ActiveWindow.ScrollColumn = 17
$ws.Range("AB10").Select()
